$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Sampling Events" sheet: remove the YM002-YM005 sampling-event rows,
#    keep only the header row and the YM001 row.
# ---------------------------------------------------------------------------
$wsEvents = $wb.Worksheets.Item("Sampling Events")
$wsEvents.Rows.Item(3).Resize(4).EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2) "Occurrences" sheet: add 7 new rows for additional species /
#    transcribed occurrences, and mark the existing + new occurrences as
#    "Present" (instead of carrying the sampling-event date range).
# ---------------------------------------------------------------------------
$wsOcc = $wb.Worksheets.Item("Occurrences")

$scientificNames = @("Bixa orellana", "Mucuna pruriens utilis", "Curcuma domestica", "Pandanus conoideus")
for ($i = 0; $i -lt 4; $i++) {
    $wsOcc.Cells.Item(7 + $i, 6).Value = $scientificNames[$i]
}

$occurrenceIds = @(
    "UNCEN-2001SS-HS004-YM001-EM006",
    "UNCEN-2001SS-HS004-YM001-EM007",
    "UNCEN-2001SS-HS004-YM001-EM008",
    "UNCEN-2001SS-HS004-YM001-EM009",
    "UNCEN-2001SS-HS004-YM001-EM010",
    "UNCEN-2001SS-HS004-YM001-EM011",
    "UNCEN-2001SS-HS004-YM001-EM012"
)
for ($i = 0; $i -lt 7; $i++) {
    $wsOcc.Cells.Item(7 + $i, 2).Value = $occurrenceIds[$i]
}

for ($r = 7; $r -le 13; $r++) {
    $wsOcc.Cells.Item($r, 1).Value  = "UNCEN-2001SS-HS004-YM001"   # A eventID
    $wsOcc.Cells.Item($r, 3).Value  = "Human Observation"          # C basisOfRecord
    $wsOcc.Cells.Item($r, 4).Value  = "2001-02-22/2001-03-22"      # D eventDate
    $wsOcc.Cells.Item($r, 5).Value  = "Plantae"                    # E kingdom
    $wsOcc.Cells.Item($r, 7).Value  = "Spesies"                    # G taxonRank
    $wsOcc.Cells.Item($r, 9).Value  = "-0.765419"                  # I decimalLatitude
    $wsOcc.Cells.Item($r, 10).Value = "133.979771"                 # J decimalLongitude
    $wsOcc.Cells.Item($r, 12).Value = "ID"                         # L countryCode
    $wsOcc.Cells.Item($r, 13).Value = "?"                          # M individualCount
}
# rows 7-10 keep the default geodeticDatum
for ($r = 7; $r -le 10; $r++) {
    $wsOcc.Cells.Item($r, 11).Value = "WGS84"                      # K geodeticDatum
}

$wsOcc.Cells.Item(11, 11).Value = "WGS85"
$wsOcc.Cells.Item(12, 11).Value = "WGS86"
$wsOcc.Cells.Item(13, 11).Value = "WGS87"

for ($r = 2; $r -le 13; $r++) {
    $wsOcc.Cells.Item($r, 16).Value = "Present"                    # P occurrenceStatus
}
